$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 716.6667
$ws.Range("J2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("N2").Value = -1226
# Row 17
$ws.Range("H17").Value = 2761.9333
$ws.Range("J17").Value = 2761.9333
$ws.Range("L17").Value = 8285.7999
$ws.Range("N17").Value = -8621.7999
# Row 19
$ws.Range("H19").Value = 2759.5
$ws.Range("J19").Value = 2453.7144
$ws.Range("L19").Value = 2453.7144
$ws.Range("N19").Value = -2803.7144
# Row 92
$ws.Range("H92").Value = 79
$ws.Range("I92").Value = 77.09090999999999
$ws.Range("J92").Value = 100
$ws.Range("K92").Value = 77.09090999999999
$ws.Range("L92").Value = 100
$ws.Range("M92").Value = 1170.90909
$ws.Range("N92").Value = -2596
# Row 96
$ws.Range("H96").Value = 377.4
$ws.Range("J96").Value = 399.5
$ws.Range("L96").Value = 1198.5
$ws.Range("N96").Value = -3944.5
# Row 98
$ws.Range("H98").Value = 1428
$ws.Range("I98").Value = 1266
$ws.Range("J98").Value = 3696
$ws.Range("K98").Value = 1266
$ws.Range("L98").Value = 3696
$ws.Range("M98").Value = 232
$ws.Range("N98").Value = -6692
# Row 100
$ws.Range("H100").Value = 1308.5333
$ws.Range("I100").Value = 1309.7273
$ws.Range("K100").Value = 1309.7273
$ws.Range("M100").Value = -768.7273
# Row 103
$ws.Range("H103").Value = 3901.3333
$ws.Range("J103").Value = 850
$ws.Range("L103").Value = 2550
$ws.Range("N103").Value = -3722
# Row 111
$ws.Range("H111").Value = 521.3333
$ws.Range("I111").Value = 469
$ws.Range("K111").Value = 1407
$ws.Range("M111").Value = 1660
# Row 122
$ws.Range("H122").Value = 1428
$ws.Range("I122").Value = 1266
$ws.Range("J122").Value = 3696
$ws.Range("K122").Value = 3798
$ws.Range("L122").Value = 11088
$ws.Range("M122").Value = -1348
$ws.Range("N122").Value = -15988
# Row 132
$ws.Range("H132").Value = 1717.5883
$ws.Range("I132").Value = 1620.2667
$ws.Range("K132").Value = 4860.800099999999
$ws.Range("M132").Value = -2330.800099999999
# Row 137
$ws.Range("H137").Value = 1064.6471
$ws.Range("I137").Value = 1000.26666
$ws.Range("K137").Value = 3000.79998
$ws.Range("M137").Value = -450.7999799999998

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 15049.833
$ws.Range("I32").Value = 13205.929
$ws.Range("J32").Value = 21503.5
$ws.Range("K32").Value = 13205.929
$ws.Range("L32").Value = 21503.5
$ws.Range("M32").Value = -12918.929
$ws.Range("N32").Value = -22077.5
# Row 61
$ws.Range("H61").Value = 5522.725
$ws.Range("I61").Value = 5629.237
$ws.Range("K61").Value = 5629.237
$ws.Range("M61").Value = -5417.237
# Row 97
$ws.Range("H97").Value = 944.86206
$ws.Range("I97").Value = 829.2083
$ws.Range("K97").Value = 829.2083
$ws.Range("M97").Value = -333.2083
# Row 125
$ws.Range("H125").Value = 24999.5
$ws.Range("J125").Value = 24999.5
$ws.Range("L125").Value = 24999.5
$ws.Range("N125").Value = -34839.5
# Row 136
$ws.Range("H136").Value = 5522.725
$ws.Range("I136").Value = 5629.237
$ws.Range("K136").Value = 16887.711
$ws.Range("M136").Value = -14337.711

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 2
$ws.Range("H2").Value = 69980
$ws.Range("J2").Value = 69980
$ws.Range("L2").Value = 69980
$ws.Range("N2").Value = -70206
# Row 26
$ws.Range("H26").Value = 20235.5
$ws.Range("I26").Value = 20235.5
$ws.Range("K26").Value = 20235.5
$ws.Range("M26").Value = -19943.5
# Row 86
$ws.Range("H86").Value = 1410.3478
$ws.Range("I86").Value = 1184.5294
$ws.Range("K86").Value = 1184.5294
$ws.Range("M86").Value = -61.5293999999999
# Row 89
$ws.Range("H89").Value = 1410.3478
$ws.Range("I89").Value = 1184.5294
$ws.Range("K89").Value = 5922.646999999999
$ws.Range("M89").Value = -306.646999999999
# Row 94
$ws.Range("H94").Value = 2042.2727
$ws.Range("I94").Value = 1244.3334
$ws.Range("J94").Value = 2999.8
$ws.Range("K94").Value = 1244.3334
$ws.Range("L94").Value = 2999.8
$ws.Range("M94").Value = -793.3334
$ws.Range("N94").Value = -3901.8
# Row 132
$ws.Range("H132").Value = 97500
$ws.Range("J132").Value = 97500
$ws.Range("L132").Value = 97500
$ws.Range("N132").Value = -107620
# Row 134
$ws.Range("H134").Value = 4921.362
$ws.Range("I134").Value = 4780.8335
$ws.Range("J134").Value = 6101.8
$ws.Range("K134").Value = 14342.5005
$ws.Range("L134").Value = 18305.4
$ws.Range("M134").Value = -11807.5005
$ws.Range("N134").Value = -23375.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 5872.364
$ws.Range("I31").Value = 4321.875
$ws.Range("J31").Value = 6758.357
$ws.Range("K31").Value = 4321.875
$ws.Range("L31").Value = 6758.357
$ws.Range("M31").Value = -4026.875
$ws.Range("N31").Value = -7348.357
# Row 34
$ws.Range("H34").Value = 5872.364
$ws.Range("I34").Value = 4321.875
$ws.Range("J34").Value = 6758.357
$ws.Range("K34").Value = 4321.875
$ws.Range("L34").Value = 6758.357
$ws.Range("M34").Value = -4119.875
$ws.Range("N34").Value = -7162.357
# Row 58
$ws.Range("H58").Value = 9410.214
$ws.Range("I58").Value = 4851.769
$ws.Range("K58").Value = 4851.769
$ws.Range("M58").Value = -4648.769
# Row 86
$ws.Range("H86").Value = 10203.053
$ws.Range("I86").Value = 10218.5
$ws.Range("J86").Value = 10159.8
$ws.Range("K86").Value = 10218.5
$ws.Range("L86").Value = 10159.8
$ws.Range("M86").Value = -9095.5
$ws.Range("N86").Value = -12405.8
# Row 89
$ws.Range("H89").Value = 10203.053
$ws.Range("I89").Value = 10218.5
$ws.Range("J89").Value = 10159.8
$ws.Range("K89").Value = 51092.5
$ws.Range("L89").Value = 50799
$ws.Range("M89").Value = -45476.5
$ws.Range("N89").Value = -62031
# Row 136
$ws.Range("H136").Value = 9410.214
$ws.Range("I136").Value = 4851.769
$ws.Range("K136").Value = 14555.307
$ws.Range("M136").Value = -12005.307
# Row 141
$ws.Range("H141").Value = 215929.73
$ws.Range("J141").Value = 215929.73
$ws.Range("L141").Value = 215929.73
$ws.Range("N141").Value = -226289.73

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 38
$ws.Range("H38").Value = 1660.3
$ws.Range("I38").Value = 1243
$ws.Range("J38").Value = 1839.1428
$ws.Range("K38").Value = 3729
$ws.Range("L38").Value = 5517.428400000001
$ws.Range("M38").Value = -3382
$ws.Range("N38").Value = -6211.428400000001
# Row 97
$ws.Range("H97").Value = 1875
$ws.Range("I97").Value = 1533
$ws.Range("J97").Value = 2331
$ws.Range("K97").Value = 4599
$ws.Range("L97").Value = 6993
$ws.Range("M97").Value = -4103
$ws.Range("N97").Value = -7985
# Row 113
$ws.Range("H113").Value = 499.33334
$ws.Range("J113").Value = 500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 51.944443
$ws.Range("J2").Value = 60.4
$ws.Range("L2").Value = 60.4
$ws.Range("N2").Value = -286.4
# Row 53
$ws.Range("H53").Value = 44039
$ws.Range("I53").Value = 44039
$ws.Range("K53").Value = 44039
$ws.Range("M53").Value = -43408

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 7632.8335
$ws.Range("J46").Value = 10000
$ws.Range("L46").Value = 10000
$ws.Range("N46").Value = -10376
# Row 93
$ws.Range("H93").Value = 13944.875
$ws.Range("I93").Value = 1292
$ws.Range("K93").Value = 1292
$ws.Range("M93").Value = -44
# Row 109
$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 7000
$ws.Range("I15").Value = 7000
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 7000
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -6712
$ws.Range("N15").ClearContents()
# Row 33
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
# Row 36
$ws.Range("H36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()
# Row 96
$ws.Range("H96").Value = 1508.8572
$ws.Range("J96").Value = 2297.3333
$ws.Range("L96").Value = 2297.3333
$ws.Range("N96").Value = -5043.3333
# Row 100
$ws.Range("H100").Value = 860.75
$ws.Range("I100").Value = 741.125
$ws.Range("J100").Value = 1100
$ws.Range("K100").Value = 1482.25
$ws.Range("L100").Value = 2200
$ws.Range("M100").Value = -941.25
$ws.Range("N100").Value = -3282
# Row 107
$ws.Range("H107").Value = 2376.4614
$ws.Range("I107").Value = 1367.1428
$ws.Range("J107").Value = 3554
$ws.Range("K107").Value = 4101.428400000001
$ws.Range("L107").Value = 10662
$ws.Range("M107").Value = -2181.428400000001
$ws.Range("N107").Value = -14502
# Row 109
$ws.Range("H109").Value = 42000
$ws.Range("J109").Value = 42000
$ws.Range("L109").Value = 42000
$ws.Range("N109").Value = -44774
